$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-10 Wednesday" "2024-07-11 Thursday"

Replace-Text "286×8=2288" "340×2=680"
Replace-Text "273×6=1638" "264×3=792"
Replace-Text "936×7=6552" "702×3=2106"
Replace-Text "356×8=2848" "767×9=6903"
Replace-Text "203×2=406" "559×8=4472"

Replace-Text "778×6=4668" "124×8=992"
Replace-Text "861×2=1722" "821×3=2463"
Replace-Text "916×4=3664" "546×5=2730"
Replace-Text "929×6=5574" "688×7=4816"
Replace-Text "543×9=4887" "233×9=2097"

Replace-Text "564×2=1128" "884×5=4420"
Replace-Text "821×4=3284" "784×2=1568"
Replace-Text "990×9=8910" "512×3=1536"
Replace-Text "748×7=5236" "702×2=1404"
Replace-Text "996×5=4980" "613×7=4291"

Replace-Text "589×6=3534" "463×2=926"
Replace-Text "241×4=964" "488×5=2440"
Replace-Text "257×3=771" "675×6=4050"
Replace-Text "125×3=375" "141×6=846"
Replace-Text "896×4=3584" "886×3=2658"

Replace-Text "399×7=2793" "413×3=1239"
Replace-Text "137×3=411" "544×6=3264"
Replace-Text "895×3=2685" "925×4=3700"
Replace-Text "990×6=5940" "317×3=951"
Replace-Text "310×3=930" "295×9=2655"
